$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "sdmx-dimension:refArea"
$ws.Range("P3").Value = "dim"
$ws.Range("P4").Value = "URI-Municipio"
